$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.039.23'
$ws.Range('E2').Value = '  +1.96%  '
$ws.Range('D3').Value = '2.312.21'
$ws.Range('E3').Value = '  +1.81%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.93'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.05'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +5.45%  '
$ws.Range('E7').Value = '  +1.75%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.514'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +3.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.17'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +5.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0797'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.67%  '
$ws.Range('E12').Value = '  +4.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '17.92'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +15.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.93'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.75%  '
$ws.Range('D15').Value = '2.684.61'
$ws.Range('E15').Value = '  +2.24%  '
$ws.Range('D16').Value = '2.297.71'
$ws.Range('E16').Value = '  +0.93%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.815'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +3.97%  '
$ws.Range('D18').Value = '42.956.43'
$ws.Range('E18').Value = '  +1.92%  '
$ws.Range('E19').Value = '  +7.53%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.16'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +2.76%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').Value = '0.0₃0905'
$ws.Range('E21').Value = '  +1.35%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.86'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.53'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.97%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.21'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +12.36%  '
$ws.Range('E25').Value = '  +0.57%  '
$ws.Range('E26').Value = '  -0.42%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.81'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +3.43%  '
$ws.Range('E28').Value = '  +0.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '167.86'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.10'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.23'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.92%  '
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.02'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.15%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.64'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.53%  '
$ws.Range('E35').Value = '  +3.61%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.09'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.84%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0694'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.86%  '
$ws.Range('E38').Value = '  +3.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.80'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +3.83%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.82'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.21%  '
$ws.Range('E41').Value = '  +0.58%  '
$ws.Range('D42').Value = '2.003.74'
$ws.Range('E42').Value = '  +2.18%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.28'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -6.94%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0287'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.45%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.23'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +7.15%  '
$ws.Range('E46').Value = '  -0.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.85'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.51'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +5.93%  '
$ws.Range('D49').Value = '2.527.73'
$ws.Range('E49').Value = '  +1.25%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.54'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +4.50%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.57'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.54%  '
